$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B50 was previously stored as text "3"; convert it to an actual number 3
$ws.Cells.Item(50, 2).Value = 3

# Add new row 51 with the new annotation data
$ws.Cells.Item(51, 1).Value = "Ruilin"

# B51 must stay as text "3" (not auto-converted to a number). Temporarily
# force a text number format so Excel keeps it as a string, then restore
# the default formatting/style so no extra styling is left on the cell.
$cellB51 = $ws.Cells.Item(51, 2)
$cellB51.NumberFormat = "@"
$cellB51.Value = "3"
$cellB51.NumberFormat = "General"
$cellB51.Style = "Normal"

$ws.Cells.Item(51, 3).Value = "无"
$ws.Cells.Item(51, 4).Value = "DIS"
$ws.Cells.Item(51, 5).Value = "WRI"
$ws.Cells.Item(51, 6).Value = "2a7301cf-d5b3-4d65-86b5-7931ca3b6163"
$ws.Cells.Item(51, 7).Value = "r1q7n9gAb_annotated.xlsx"
$ws.Cells.Item(51, 8).Value = "It would be beneficial for the clarity if authors define what they mean by convergence (normalised weight vector, angle, whichever path seems most natural) as early in the paper as possible."
